# Ticket 86 - Move, remove, and copy existing Excel hyperlinks.
# Adds a new "Shift" worksheet (after the existing sheets) that demonstrates
# JETT's ability to shift and copy hyperlinks along with a <jt:for> loop.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet, so it ends up last in tab order.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Shift"

# Header row.
$ws.Range("A1").Value = "Shift!"
$ws.Range("B1").Value = "Copy!"

# Loop tag row, with a hyperlink ("Example") that will be shifted/copied by the loop.
$ws.Range("A2").Value = '<jt:for var="i" start="1" end="10">'
$ws.Range("B2").Value = "Example"
$ws.Hyperlinks.Add($ws.Range("B2"), "http://www.jett.com")
$ws.Range("C2").Value = "</jt:for>"

# Closing row with another hyperlink ("JETT") outside the loop.
$ws.Range("B3").Value = "JETT"
$ws.Hyperlinks.Add($ws.Range("B3"), "http://www.jett.com")
